$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-11-25 Monday" "2024-11-26 Tuesday"

Replace-Text "49÷5=9, 4" "76÷5=15, 1"
Replace-Text "64÷7=9, 1" "95÷8=11, 7"
Replace-Text "63÷5=12, 3" "23÷2=11, 1"
Replace-Text "61÷8=7, 5" "92÷5=18, 2"
Replace-Text "84÷7=12, 0" "53÷4=13, 1"
Replace-Text "97÷4=24, 1" "36÷6=6, 0"
Replace-Text "79÷3=26, 1" "13÷3=4, 1"
Replace-Text "77÷8=9, 5" "42÷3=14, 0"
Replace-Text "21÷9=2, 3" "75÷5=15, 0"
Replace-Text "57÷7=8, 1" "85÷3=28, 1"
Replace-Text "84÷2=42, 0" "60÷2=30, 0"
Replace-Text "90÷8=11, 2" "72÷5=14, 2"
Replace-Text "48÷6=8, 0" "46÷2=23, 0"
Replace-Text "32÷2=16, 0" "56÷7=8, 0"
Replace-Text "98÷5=19, 3" "31÷4=7, 3"
Replace-Text "19÷4=4, 3" "49÷3=16, 1"
Replace-Text "41÷3=13, 2" "95÷4=23, 3"
Replace-Text "65÷5=13, 0" "42÷9=4, 6"
Replace-Text "30÷4=7, 2" "57÷9=6, 3"
Replace-Text "43÷2=21, 1" "94÷7=13, 3"
Replace-Text "34÷2=17, 0" "58÷2=29, 0"
Replace-Text "28÷4=7, 0" "20÷8=2, 4"
Replace-Text "53÷9=5, 8" "74÷6=12, 2"
Replace-Text "54÷6=9, 0" "20÷8=2, 4"
Replace-Text "22÷7=3, 1" "42÷9=4, 6"
